$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.039748046762647
$ws.Cells.Item(2, 4).Value = 1.047694317005407
$ws.Cells.Item(2, 5).Value = 1.05303349684609
$ws.Cells.Item(2, 6).Value = 1.059851346094653
$ws.Cells.Item(2, 9).Value = 1.03834399213758
$ws.Cells.Item(2, 10).Value = 1.044838385046212
$ws.Cells.Item(2, 11).Value = 1.050456238647997
$ws.Cells.Item(2, 12).Value = 1.055780583095578
$ws.Cells.Item(2, 13).Value = 1.062579723669453
$ws.Cells.Item(2, 14).Value = 1.018852173824719
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.040726314777273
$ws.Cells.Item(3, 4).Value = 1.048472650231064
$ws.Cells.Item(3, 5).Value = 1.05397132229688
$ws.Cells.Item(3, 6).Value = 1.060807409537149
$ws.Cells.Item(3, 9).Value = 1.03853798902463
$ws.Cells.Item(3, 10).Value = 1.045461842691475
$ws.Cells.Item(3, 11).Value = 1.051046554627963
$ws.Cells.Item(3, 12).Value = 1.056531052649202
$ws.Cells.Item(3, 13).Value = 1.063349739224914
$ws.Cells.Item(3, 14).Value = 1.019060091420595
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.041359608897503
$ws.Cells.Item(4, 4).Value = 1.04897611339514
$ws.Cells.Item(4, 5).Value = 1.054578732856905
$ws.Cells.Item(4, 6).Value = 1.061426493631996
$ws.Cells.Item(4, 9).Value = 1.038661701488782
$ws.Cells.Item(4, 10).Value = 1.045864941083136
$ws.Cells.Item(4, 11).Value = 1.051427699749322
$ws.Cells.Item(4, 12).Value = 1.057016604298479
$ws.Cells.Item(4, 13).Value = 1.06384781341129
$ws.Cells.Item(4, 14).Value = 1.019194478014975
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.041625914287945
$ws.Cells.Item(5, 4).Value = 1.049187727247312
$ws.Cells.Item(5, 5).Value = 1.054834224734483
$ws.Cells.Item(5, 6).Value = 1.061686862445847
$ws.Cells.Item(5, 9).Value = 1.038713274691157
$ws.Cells.Item(5, 10).Value = 1.046034326125016
$ws.Cells.Item(5, 11).Value = 1.051587733762158
$ws.Cells.Item(5, 12).Value = 1.057220716717538
$ws.Cells.Item(5, 13).Value = 1.064057160293259
$ws.Cells.Item(5, 14).Value = 1.019250937912808
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.041670632133337
$ws.Cells.Item(6, 4).Value = 1.049223255641274
$ws.Cells.Item(6, 5).Value = 1.054877130915304
$ws.Cells.Item(6, 6).Value = 1.061730585688692
$ws.Cells.Item(6, 9).Value = 1.038721908500968
$ws.Cells.Item(6, 10).Value = 1.046062762049092
$ws.Cells.Item(6, 11).Value = 1.051614592455765
$ws.Cells.Item(6, 12).Value = 1.057254987271404
$ws.Cells.Item(6, 13).Value = 1.064092307974283
$ws.Cells.Item(6, 14).Value = 1.019260415644495
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.041363167014737
$ws.Cells.Item(7, 4).Value = 1.048978941156647
$ws.Cells.Item(7, 5).Value = 1.054582146217032
$ws.Cells.Item(7, 6).Value = 1.061429972278362
$ws.Cells.Item(7, 9).Value = 1.038662392324542
$ws.Cells.Item(7, 10).Value = 1.045867204719999
$ws.Cells.Item(7, 11).Value = 1.051429838916371
$ws.Cells.Item(7, 12).Value = 1.057019331712856
$ws.Cells.Item(7, 13).Value = 1.063850610886842
$ws.Cells.Item(7, 14).Value = 1.01919523257764
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.040078596686323
$ws.Cells.Item(8, 4).Value = 1.047957392912844
$ws.Cells.Item(8, 5).Value = 1.053350319659595
$ws.Cells.Item(8, 6).Value = 1.060174358924285
$ws.Cells.Item(8, 9).Value = 1.038409930024396
$ws.Cells.Item(8, 10).Value = 1.045049151324848
$ws.Cells.Item(8, 11).Value = 1.050655909715233
$ws.Cells.Item(8, 12).Value = 1.056034218051828
$ws.Cells.Item(8, 13).Value = 1.062839990130493
$ws.Cells.Item(8, 14).Value = 1.01892247141613
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.037817266637019
$ws.Cells.Item(9, 4).Value = 1.046156043444106
$ws.Cells.Item(9, 5).Value = 1.051184135404372
$ws.Cells.Item(9, 6).Value = 1.057965283892495
$ws.Cells.Item(9, 9).Value = 1.037951175198821
$ws.Cells.Item(9, 10).Value = 1.043605218730754
$ws.Cells.Item(9, 11).Value = 1.049285842981846
$ws.Cells.Item(9, 12).Value = 1.054297963381674
$ws.Cells.Item(9, 13).Value = 1.061057834522831
$ws.Cells.Item(9, 14).Value = 1.018440696557723
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.036311264091208
$ws.Cells.Item(10, 4).Value = 1.044954374980873
$ws.Cells.Item(10, 5).Value = 1.049743070447119
$ws.Cells.Item(10, 6).Value = 1.056494972367548
$ws.Cells.Item(10, 9).Value = 1.037636040890978
$ws.Cells.Item(10, 10).Value = 1.042641018239332
$ws.Cells.Item(10, 11).Value = 1.048368282763213
$ws.Cells.Item(10, 12).Value = 1.053140274708251
$ws.Cells.Item(10, 13).Value = 1.059868908233322
$ws.Cells.Item(10, 14).Value = 1.018118768657838
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.035659523956296
$ws.Cells.Item(11, 4).Value = 1.044433871671409
$ws.Cells.Item(11, 5).Value = 1.049119811872187
$ws.Cells.Item(11, 6).Value = 1.055858895510704
$ws.Cells.Item(11, 9).Value = 1.037497386456943
$ws.Cells.Item(11, 10).Value = 1.042223143784029
$ws.Cells.Item(11, 11).Value = 1.047969988492163
$ws.Cells.Item(11, 12).Value = 1.052638949941539
$ws.Cells.Item(11, 13).Value = 1.059353907081846
$ws.Cells.Item(11, 14).Value = 1.017979197428263
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.035417494618933
$ws.Cells.Item(12, 4).Value = 1.044240508718657
$ws.Cells.Item(12, 5).Value = 1.048888416876817
$ws.Cells.Item(12, 6).Value = 1.055622716094655
$ws.Cells.Item(12, 9).Value = 1.037445554173269
$ws.Cells.Item(12, 10).Value = 1.042067871833147
$ws.Cells.Item(12, 11).Value = 1.047821897195548
$ws.Cells.Item(12, 12).Value = 1.052452730750161
$ws.Cells.Item(12, 13).Value = 1.059162585364358
$ws.Cells.Item(12, 14).Value = 1.017927328551582
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.035469408192008
$ws.Cells.Item(13, 4).Value = 1.044281986834707
$ws.Cells.Item(13, 5).Value = 1.048938046850164
$ws.Cells.Item(13, 6).Value = 1.055673373388888
$ws.Cells.Item(13, 9).Value = 1.037456687291806
$ws.Cells.Item(13, 10).Value = 1.042101180660921
$ws.Cells.Item(13, 11).Value = 1.047853669922846
$ws.Cells.Item(13, 12).Value = 1.052492675604819
$ws.Cells.Item(13, 13).Value = 1.059203625745572
$ws.Cells.Item(13, 14).Value = 1.01793845576955
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.035639516590902
$ws.Cells.Item(14, 4).Value = 1.044417888729777
$ws.Cells.Item(14, 5).Value = 1.049100682419717
$ws.Cells.Item(14, 6).Value = 1.055839371048316
$ws.Cells.Item(14, 9).Value = 1.037493108709935
$ws.Cells.Item(14, 10).Value = 1.042210310069386
$ws.Cells.Item(14, 11).Value = 1.047957750214615
$ws.Cells.Item(14, 12).Value = 1.052623557099168
$ws.Cells.Item(14, 13).Value = 1.059338092921818
$ws.Cells.Item(14, 14).Value = 1.017974910457613
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.035744333422255
$ws.Cells.Item(15, 4).Value = 1.044501619113832
$ws.Cells.Item(15, 5).Value = 1.049200902299059
$ws.Cells.Item(15, 6).Value = 1.055941659357266
$ws.Cells.Item(15, 9).Value = 1.037515505455379
$ws.Cells.Item(15, 10).Value = 1.042277541062932
$ws.Cells.Item(15, 11).Value = 1.04802185805489
$ws.Cells.Item(15, 12).Value = 1.052704196889892
$ws.Cells.Item(15, 13).Value = 1.059420938989119
$ws.Cells.Item(15, 14).Value = 1.017997367969513
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.036354525756309
$ws.Cells.Item(16, 4).Value = 1.044988915524613
$ws.Cells.Item(16, 5).Value = 1.049784449575712
$ws.Cells.Item(16, 6).Value = 1.056537198955745
$ws.Cells.Item(16, 9).Value = 1.037645196641679
$ws.Cells.Item(16, 10).Value = 1.042668743470866
$ws.Cells.Item(16, 11).Value = 1.048394695570563
$ws.Cells.Item(16, 12).Value = 1.053173545255156
$ws.Cells.Item(16, 13).Value = 1.059903083276516
$ws.Cells.Item(16, 14).Value = 1.018128027892014
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.036737382436612
$ws.Cells.Item(17, 4).Value = 1.045294538191271
$ws.Cells.Item(17, 5).Value = 1.050150690070057
$ws.Cells.Item(17, 6).Value = 1.056910920631883
$ws.Cells.Item(17, 9).Value = 1.037725960098955
$ws.Cells.Item(17, 10).Value = 1.042914036039829
$ws.Cells.Item(17, 11).Value = 1.04862830368679
$ws.Cells.Item(17, 12).Value = 1.053467945528465
$ws.Cells.Item(17, 13).Value = 1.060205469907378
$ws.Cells.Item(17, 14).Value = 1.018209940937706
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.036960731818386
$ws.Cells.Item(18, 4).Value = 1.04547278596061
$ws.Cells.Item(18, 5).Value = 1.050364382520109
$ws.Cells.Item(18, 6).Value = 1.057128961835547
$ws.Cells.Item(18, 9).Value = 1.037772855829567
$ws.Cells.Item(18, 10).Value = 1.043057075367828
$ws.Cells.Item(18, 11).Value = 1.048764468287558
$ws.Cells.Item(18, 12).Value = 1.05363966051054
$ws.Cells.Item(18, 13).Value = 1.060381828691041
$ws.Cells.Item(18, 14).Value = 1.018257702592715
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.037036894200354
$ws.Cells.Item(19, 4).Value = 1.045533560966029
$ws.Cells.Item(19, 5).Value = 1.050437258061835
$ws.Cells.Item(19, 6).Value = 1.057203317662339
$ws.Cells.Item(19, 9).Value = 1.037788810047495
$ws.Cells.Item(19, 10).Value = 1.043105842009122
$ws.Cells.Item(19, 11).Value = 1.048810880747341
$ws.Cells.Item(19, 12).Value = 1.053698210225714
$ws.Cells.Item(19, 13).Value = 1.06044195935995
$ws.Cells.Item(19, 14).Value = 1.018273985218053
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.036696301856237
$ws.Cells.Item(20, 4).Value = 1.04526174949698
$ws.Cells.Item(20, 5).Value = 1.050111388624773
$ws.Cells.Item(20, 6).Value = 1.056870818066729
$ws.Cells.Item(20, 9).Value = 1.037717316891999
$ws.Cells.Item(20, 10).Value = 1.042887722153721
$ws.Cells.Item(20, 11).Value = 1.048603249573095
$ws.Cells.Item(20, 12).Value = 1.053436359543546
$ws.Cells.Item(20, 13).Value = 1.060173028561741
$ws.Cells.Item(20, 14).Value = 1.018201154182765
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.035589422357923
$ws.Cells.Item(21, 4).Value = 1.044377869685803
$ws.Cells.Item(21, 5).Value = 1.049052787226138
$ws.Cells.Item(21, 6).Value = 1.055790486441482
$ws.Cells.Item(21, 9).Value = 1.03748239261123
$ws.Cells.Item(21, 10).Value = 1.042178175684885
$ws.Cells.Item(21, 11).Value = 1.04792710521579
$ws.Cells.Item(21, 12).Value = 1.052585015888354
$ws.Cells.Item(21, 13).Value = 1.059298496437811
$ws.Cells.Item(21, 14).Value = 1.017964176173499
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.034893807365584
$ws.Cells.Item(22, 4).Value = 1.043821996160299
$ws.Cells.Item(22, 5).Value = 1.048387845160625
$ws.Cells.Item(22, 6).Value = 1.055111747861009
$ws.Cells.Item(22, 9).Value = 1.03733277823287
$ws.Cells.Item(22, 10).Value = 1.041731739110944
$ws.Cells.Item(22, 11).Value = 1.04750113594184
$ws.Cells.Item(22, 12).Value = 1.052049714690507
$ws.Cells.Item(22, 13).Value = 1.058748485540699
$ws.Cells.Item(22, 14).Value = 1.017815028939213
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.035262535133791
$ws.Cells.Item(23, 4).Value = 1.044116688431113
$ws.Cells.Item(23, 5).Value = 1.048740282279773
$ws.Cells.Item(23, 6).Value = 1.055471511447007
$ws.Cells.Item(23, 9).Value = 1.037412272333851
$ws.Cells.Item(23, 10).Value = 1.041968433355004
$ws.Cells.Item(23, 11).Value = 1.047727030643019
$ws.Cells.Item(23, 12).Value = 1.052333490322684
$ws.Cells.Item(23, 13).Value = 1.059040071379413
$ws.Cells.Item(23, 14).Value = 1.017894108810926
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.036714864292313
$ws.Cells.Item(24, 4).Value = 1.045276565348418
$ws.Cells.Item(24, 5).Value = 1.050129147038198
$ws.Cells.Item(24, 6).Value = 1.056888938517758
$ws.Cells.Item(24, 9).Value = 1.037721223040948
$ws.Cells.Item(24, 10).Value = 1.042899612376443
$ws.Cells.Item(24, 11).Value = 1.04861457074239
$ws.Cells.Item(24, 12).Value = 1.053450631902043
$ws.Cells.Item(24, 13).Value = 1.060187687465754
$ws.Cells.Item(24, 14).Value = 1.018205124591144
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.038401604247778
$ws.Cells.Item(25, 4).Value = 1.046621876234527
$ws.Cells.Item(25, 5).Value = 1.051743611725325
$ws.Cells.Item(25, 6).Value = 1.058535963776689
$ws.Cells.Item(25, 9).Value = 1.038071415808021
$ws.Cells.Item(25, 10).Value = 1.043978791610808
$ws.Cells.Item(25, 11).Value = 1.049640779362278
$ws.Cells.Item(25, 12).Value = 1.05474686358773
$ws.Cells.Item(25, 13).Value = 1.061518713825345
$ws.Cells.Item(25, 14).Value = 1.018565379496608
